# Region XII_HEALTH.xlsx edit: refresh "Status as of" date, add header
# formatting (centered/wrapped bold headers, yellow highlight on the new
# status column), set explicit column widths, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Update the "Status as of ..." header text (AA1)
# ------------------------------------------------------------------
$ws.Range("AA1").Value = "Status as of July 11, 2025"

# ------------------------------------------------------------------
# 2) Explicit column widths (A:AA) -- ColumnWidth is offset from the
#    stored OOXML <col width> by the default 0.8333... padding, so we
#    subtract it here to land on the exact target widths.
# ------------------------------------------------------------------
$pad = 0.8333333333333333
$ws.Columns.Item(1).ColumnWidth = 31 - $pad
$ws.Columns.Item(2).ColumnWidth = 12 - $pad
$ws.Columns.Item(3).ColumnWidth = 21 - $pad
$ws.Columns.Item(4).ColumnWidth = 11 - $pad
$ws.Columns.Item(5).ColumnWidth = 34 - $pad
$ws.Columns.Item(6).ColumnWidth = 33 - $pad
$ws.Columns.Item(7).ColumnWidth = 6 - $pad
$ws.Columns.Item(8).ColumnWidth = 20 - $pad
$ws.Columns.Item(9).ColumnWidth = 23 - $pad
$ws.Columns.Item(10).ColumnWidth = 42 - $pad
$ws.Columns.Item(11).ColumnWidth = 20 - $pad
$ws.Columns.Item(12).ColumnWidth = 7 - $pad
$ws.Columns.Item(13).ColumnWidth = 23 - $pad
$ws.Columns.Item(14).ColumnWidth = 11 - $pad
$ws.Columns.Item(15).ColumnWidth = 26 - $pad
$ws.Columns.Item(16).ColumnWidth = 26 - $pad
$ws.Columns.Item(17).ColumnWidth = 27 - $pad
$ws.Columns.Item(18).ColumnWidth = 15 - $pad
$ws.Columns.Item(19).ColumnWidth = 14 - $pad
$ws.Columns.Item(20).ColumnWidth = 31 - $pad
$ws.Columns.Item(21).ColumnWidth = 27 - $pad
$ws.Columns.Item(22).ColumnWidth = 21 - $pad
$ws.Columns.Item(23).ColumnWidth = 33 - $pad
$ws.Columns.Item(24).ColumnWidth = 31 - $pad
$ws.Columns.Item(25).ColumnWidth = 36 - $pad
$ws.Columns.Item(26).ColumnWidth = 47 - $pad
$ws.Columns.Item(27).ColumnWidth = 28 - $pad

# ------------------------------------------------------------------
# 3) Header row (A1:AA1): center + middle aligned, wrapped text
# ------------------------------------------------------------------
$headerRange = $ws.Range("A1:AA1")
$headerRange.WrapText = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# ------------------------------------------------------------------
# 4) Yellow highlight fill for the new "Status as of" column (AA1:AA8)
# ------------------------------------------------------------------
$ws.Range("AA1:AA8").Interior.Color = 65535

# ------------------------------------------------------------------
# 5) Freeze the header row (split below row 1), keep A1 selected
# ------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Host "edit complete"
